# Updates cryptos list values (price / 1h volume %) to match the
# latest scrape; a couple of rows (PancakeSwap / EnergySwap) also
# swapped rank position, so their Coin/Link/Price/Volume cells move too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "0.9558"); Excel
# would silently convert such a string to a real number on assignment,
# so those writes get a leading apostrophe (the same quote-prefix a
# user gets by typing '0.9558 into a cell) to keep it stored as text,
# exactly like the source cell already was. Values such as "20.537.75"
# or "1.474.58" have more than one dot and are never read as numbers,
# so they are assigned as plain strings.

$ws.Range("D2").Value = '20.537.75'
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").Value = '1.474.58'
$ws.Range("E3").Value = '  +2.26%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = "'" + '0.9558'
$ws.Range("E5").Value = '  +4.17%  '

$ws.Range("D6").Value = "'" + '278.04'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").Value = "'" + '0.3620'
$ws.Range("E7").Value = '  -1.42%  '

$ws.Range("D8").Value = "'" + '0.3071'
$ws.Range("E8").Value = '  -2.04%  '

$ws.Range("D9").Value = "'" + '39.68'
$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("D10").Value = "'" + '1.066'
$ws.Range("E10").Value = '  +3.92%  '

$ws.Range("D11").Value = "'" + '0.06662'
$ws.Range("E11").Value = '  +1.88%  '

$ws.Range("D12").Value = "'" + '1.002'
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").Value = "'" + '5.533'
$ws.Range("E13").Value = '  +2.23%  '

$ws.Range("E14").Value = '  +2.94%  '

$ws.Range("D15").Value = "'" + '6.200'
$ws.Range("E15").Value = '  +2.11%  '

$ws.Range("D16").Value = "'" + '0.9553'
$ws.Range("E16").Value = '  +2.16%  '

$ws.Range("E17").Value = '  +0.99%  '

$ws.Range("D18").Value = '1.475.51'
$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("D19").Value = "'" + '0.05923'
$ws.Range("E19").Value = '  +5.03%  '

$ws.Range("D20").Value = "'" + '69.12'
$ws.Range("E20").Value = '  +2.03%  '

$ws.Range("D21").Value = "'" + '5.513'
$ws.Range("E21").Value = '  +1.48%  '

$ws.Range("D22").Value = "'" + '14.53'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("E23").Value = '  +2.85%  '

$ws.Range("D24").Value = "'" + '2.249'
$ws.Range("E24").Value = '  -1.16%  '

$ws.Range("D25").Value = '20.557.17'
$ws.Range("E25").Value = '  +1.68%  '

$ws.Range("D26").Value = "'" + '143.42'
$ws.Range("E26").Value = '  +5.13%  '

$ws.Range("E27").Value = '  -2.66%  '

$ws.Range("D28").Value = "'" + '17.18'
$ws.Range("E28").Value = '  +1.17%  '

$ws.Range("D29").Value = '1.636.66'
$ws.Range("E29").Value = '  +2.19%  '

$ws.Range("D30").Value = "'" + '113.82'
$ws.Range("E30").Value = '  +2.54%  '

$ws.Range("D31").Value = "'" + '3.924'
$ws.Range("E31").Value = '  +3.99%  '

$ws.Range("D32").Value = "'" + '4.983'
$ws.Range("E32").Value = '  +2.83%  '

$ws.Range("D33").Value = "'" + '0.8100'
$ws.Range("E33").Value = '  -0.70%  '

$ws.Range("D34").Value = "'" + '0.08001'
$ws.Range("E34").Value = '  +3.96%  '

$ws.Range("D35").Value = "'" + '1.515'
$ws.Range("E35").Value = '  +2.08%  '

$ws.Range("D36").Value = "'" + '1.218'
$ws.Range("E36").Value = '  +7.62%  '

$ws.Range("D37").Value = "'" + '0.05799'
$ws.Range("E37").Value = '  -3.76%  '

$ws.Range("D38").Value = "'" + '4.740'
$ws.Range("E38").Value = '  +0.53%  '

$ws.Range("D39").Value = "'" + '0.02058'
$ws.Range("E39").Value = '  +3.24%  '

$ws.Range("D40").Value = "'" + '10.39'
$ws.Range("E40").Value = '  +1.42%  '

$ws.Range("D41").Value = "'" + '0.9564'
$ws.Range("E41").Value = '  +2.23%  '

$ws.Range("D42").Value = "'" + '0.1882'
$ws.Range("E42").Value = '  +2.75%  '

$ws.Range("D43").Value = "'" + '7.431'
$ws.Range("E43").Value = '  +3.00%  '

$ws.Range("D44").Value = "'" + '0.5289'
$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'" + '12.28'
$ws.Range("E45").Value = '  +1.76%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'" + '3.522'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").Value = "'" + '118.43'
$ws.Range("E47").Value = '  -0.82%  '

$ws.Range("D48").Value = "'" + '0.5208'
$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("D49").Value = "'" + '1.816'
$ws.Range("E49").Value = '  +2.39%  '

$ws.Range("E50").Value = '  +2.28%  '

$ws.Range("D51").Value = "'" + '0.9857'
$ws.Range("E51").Value = '  -0.90%  '
